$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

Set-TextValue $ws.Range("D2") '68.007.77'
Set-TextValue $ws.Range("E2") '  +1.23%  '

Set-TextValue $ws.Range("D3") '3.527.45'
Set-TextValue $ws.Range("E3") '  +0.20%  '

Set-TextValue $ws.Range("E4") '  +0.00%  '

Set-TextValue $ws.Range("D5") '601.07'
Set-TextValue $ws.Range("E5") '  +1.25%  '

Set-TextValue $ws.Range("D6") '183.98'
Set-TextValue $ws.Range("E6") '  +5.74%  '

Set-TextValue $ws.Range("E7") '  +0.06%  '

Set-TextValue $ws.Range("E9") '  +4.34%  '

Set-TextValue $ws.Range("D10") '7.15'
Set-TextValue $ws.Range("E10") '  -1.63%  '

Set-TextValue $ws.Range("D11") '0.446'
Set-TextValue $ws.Range("E11") '  +2.06%  '

Set-TextValue $ws.Range("D12") '4.141.03'
Set-TextValue $ws.Range("E12") '  +0.31%  '

Set-TextValue $ws.Range("D13") '32.72'
Set-TextValue $ws.Range("E13") '  +12.43%  '

Set-TextValue $ws.Range("E14") '  -0.28%  '

Set-TextValue $ws.Range("D15") '67.993.28'
Set-TextValue $ws.Range("E15") '  +1.23%  '

Set-TextValue $ws.Range("E16") '  +0.56%  '

Set-TextValue $ws.Range("D17") '3.535.46'
Set-TextValue $ws.Range("E17") '  +1.58%  '

Set-TextValue $ws.Range("E18") '  +1.55%  '

Set-TextValue $ws.Range("D19") '14.85'
Set-TextValue $ws.Range("E19") '  +4.12%  '

Set-TextValue $ws.Range("D20") '400.24'
Set-TextValue $ws.Range("E20") '  +1.36%  '

Set-TextValue $ws.Range("D21") '8.12'
Set-TextValue $ws.Range("E21") '  +1.47%  '

Set-TextValue $ws.Range("D22") '73.93'
Set-TextValue $ws.Range("E22") '  +1.11%  '

Set-TextValue $ws.Range("E23") '  +1.14%  '

Set-TextValue $ws.Range("D24") '0.999'
Set-TextValue $ws.Range("E24") '  -0.31%  '

$ws.Range("B25").Value = 'LEO'
$ws.Range("C25").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range("D25") '5.69'
Set-TextValue $ws.Range("E25") '  +0.22%  '

$ws.Range("B26").Value = 'PEPE'
$ws.Range("C26").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws.Range("D26") '0.0000125'
Set-TextValue $ws.Range("E26") '  +2.21%  '

Set-TextValue $ws.Range("D27") '10.65'
Set-TextValue $ws.Range("E27") '  +3.43%  '

Set-TextValue $ws.Range("E28") '  -1.08%  '

Set-TextValue $ws.Range("E29") '  -0.16%  '

$ws.Range("B30").Value = 'Fetch.AI'
$ws.Range("C30").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range("D30") '1.49'
Set-TextValue $ws.Range("E30") '  +1.92%  '

$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D31") '6.31'
Set-TextValue $ws.Range("E31") '  +0.42%  '

Set-TextValue $ws.Range("E32") '  +1.11%  '

Set-TextValue $ws.Range("D33") '24.19'
Set-TextValue $ws.Range("E33") '  +1.14%  '

Set-TextValue $ws.Range("D34") '7.47'
Set-TextValue $ws.Range("E34") '  +1.45%  '

Set-TextValue $ws.Range("E35") '  -0.08%  '

Set-TextValue $ws.Range("E36") '  +1.95%  '

Set-TextValue $ws.Range("D37") '164.02'
Set-TextValue $ws.Range("E37") '  +0.66%  '

Set-TextValue $ws.Range("D38") '0.883'
Set-TextValue $ws.Range("E38") '  -1.74%  '

Set-TextValue $ws.Range("E39") '  +2.27%  '

Set-TextValue $ws.Range("E40") '  +3.87%  '

Set-TextValue $ws.Range("D41") '2.82'
Set-TextValue $ws.Range("E41") '  +7.16%  '

$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D42") '4.77'
Set-TextValue $ws.Range("E42") '  +1.58%  '

$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D43") '27.27'
Set-TextValue $ws.Range("E43") '  +3.03%  '

Set-TextValue $ws.Range("D44") '2.902.55'
Set-TextValue $ws.Range("E44") '  +3.39%  '

Set-TextValue $ws.Range("D45") '27.56'
Set-TextValue $ws.Range("E45") '  +0.10%  '

Set-TextValue $ws.Range("E46") '  -0.22%  '

Set-TextValue $ws.Range("D47") '42.59'
Set-TextValue $ws.Range("E47") '  -0.87%  '

Set-TextValue $ws.Range("D48") '351.94'
Set-TextValue $ws.Range("E48") '  +4.32%  '

Set-TextValue $ws.Range("D49") '0.0306'
Set-TextValue $ws.Range("E49") '  +0.43%  '

Set-TextValue $ws.Range("E50") '  -0.41%  '

Set-TextValue $ws.Range("D51") '33.65'
Set-TextValue $ws.Range("E51") '  +0.08%  '
